# Leave 3/9/2023 12:08 AM
# Applies the leave-card update: fills in Dec 2022 / Jan-Feb 2023 leave entries,
# and extends the yearly table with additional blank month rows through Oct 2026.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# ---------------------------------------------------------------------------
# 1) Extend the table by 3 rows (old last row 141 -> new last row 144).
#    First copy formats (doesn't require the destination to be populated
#    beforehand), then resize the table, then fix up the formulas.
# ---------------------------------------------------------------------------

# Push the "totals" style row (was row 141) down to the new last row (144)
$ws.Range("A141:K141").Copy()
$ws.Range("A144:K144").PasteSpecial(-4122)

# Fill rows 141-143 with the standard blank data-row formatting (copy of row 140)
$ws.Range("A140:K140").Copy()
$ws.Range("A141:K143").PasteSpecial(-4122)

# Restore the calculated-column formula text on column G for the new/moved rows
foreach ($r in 141,142,143,144) {
    $ws.Range("G$r").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
}

# Grow the table reference to include the 3 new rows
$lo.Resize($ws.Range("A8:K144"))

# ---------------------------------------------------------------------------
# 2) Fill in the actual leave data for Dec 2022 / Jan-Feb 2023
# ---------------------------------------------------------------------------

# Row 94: continuation of Dec 2022 - Sick Leave (1 day), approved 12/28/2022
$ws.Range("B94").Value = "SL(1-0-0)"
$ws.Range("H94").Value = 1
$ws.Range("K90").Copy()
$ws.Range("K94").PasteSpecial(-4122)
$ws.Range("K94").Value = 44923

# Row 95: "2023" year header
$ws.Range("A76").Copy()
$ws.Range("A95").PasteSpecial(-4122)
$ws.Range("A95").Value = "'2023"

# Row 96: Jan 2023 - Sick Leave (3 days)
$ws.Range("A96").Value = 44927
$ws.Range("B96").Value = "SL(3-0-0)"
$ws.Range("C96").Value = 1.25
$ws.Range("H96").Value = 3
$ws.Range("K96").Value = "1/9,10,11"

# Row 97: Feb 2023 - Special Privilege Leave (1 day), approved 2/20/2023
$ws.Range("A97").Value = 44958
$ws.Range("B97").Value = "SP(1-0-0)"
$ws.Range("C97").Value = 1.25
$ws.Range("K84").Copy()
$ws.Range("K97").PasteSpecial(-4122)
$ws.Range("K97").Value = 44977

# Row 98: continuation - Special Privilege Leave (2 days), period 2/8,9/2023
$ws.Range("B98").Value = "SP(2-0-0)"
$ws.Range("K84").Copy()
$ws.Range("K98").PasteSpecial(-4122)
$ws.Range("K98").Value = "2/8,9/2023"

# ---------------------------------------------------------------------------
# 3) Fill in the monthly PERIOD dates down column A through the new rows
# ---------------------------------------------------------------------------

$periodDates = @{
    99  = 44986
    100 = 45017
    101 = 45047
    102 = 45078
    103 = 45108
    104 = 45139
    105 = 45170
    106 = 45200
    107 = 45231
    108 = 45261
    109 = 45292
    110 = 45323
    111 = 45352
    112 = 45383
    113 = 45413
    114 = 45444
    115 = 45474
    116 = 45505
    117 = 45536
    118 = 45566
    119 = 45597
    120 = 45627
    121 = 45658
    122 = 45689
    123 = 45717
    124 = 45748
    125 = 45778
    126 = 45809
    127 = 45839
    128 = 45870
    129 = 45901
    130 = 45931
    131 = 45962
    132 = 45992
    133 = 46023
    134 = 46054
    135 = 46082
    136 = 46113
    137 = 46143
    138 = 46174
    139 = 46204
    140 = 46235
    141 = 46266
    142 = 46296
}

foreach ($r in $periodDates.Keys) {
    $ws.Range("A$r").Value = $periodDates[$r]
}

# ---------------------------------------------------------------------------
# 4) Recalculate and restore the previous selection
# ---------------------------------------------------------------------------

$excel.Calculate()
$ws.Range("B99").Select()
